# Actualización automática 2026-01-07 11:30:09
# Inserts a new client "CASTILLO CHOEZ CRISTIAN MARIANO" (for asesor
# OFICINA-CATAECSA) as row 18 in both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets, pushing every subsequent client row down by one,
# and refreshes the trailing "0 de NN" / "N de NN" summary-row captions on
# "VENTAS POR GRUPO" to reflect the new total client count (62 -> 63).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Push existing row 18 ("COBO FOLLECO JORGE ERNESTO"...) and everything
# below it down by one row.
$ws1.Rows.Item(18).Insert()

# Populate the freshly inserted row.
$ws1.Cells.Item(18, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(18, 2).Value = "CASTILLO CHOEZ CRISTIAN MARIANO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(18, $col).Value = 0
}

# The trailing summary row (formerly row 64, now row 65) reports counts
# like "0 de 62" / "1 de 62" -- bump the denominator to the new total of
# 63 client rows.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(65, $col)
    $oldText = $cell.Value()
    $newText = $oldText.Replace("de 62", "de 63")
    $cell.Value = $newText
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(18).Insert()

$ws2.Cells.Item(18, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(18, 2).Value = "CASTILLO CHOEZ CRISTIAN MARIANO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(18, $col).Value = 0
}
